$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheetId 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4506
$ws1.Range("F3").Value = 857
$ws1.Range("F5").Value = 38
$ws1.Range("F7").Value = 153
$ws1.Range("G7").Value = 45
$ws1.Range("F8").Value = 633
$ws1.Range("F9").Value = 24
$ws1.Range("F10").Value = 195
$ws1.Range("F11").Value = 1341
$ws1.Range("F13").Value = 2957
$ws1.Range("F14").Value = 447
$ws1.Range("F15").Value = 659

# --- Sheet "全部类型" (sheetId 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4506
$ws4.Range("F3").Value = 857
$ws4.Range("F5").Value = 38
$ws4.Range("F7").Value = 153
$ws4.Range("G7").Value = 45
$ws4.Range("F8").Value = 633
$ws4.Range("F9").Value = 24
$ws4.Range("F11").Value = 195
$ws4.Range("F12").Value = 1341
$ws4.Range("F14").Value = 2957
$ws4.Range("F15").Value = 447
$ws4.Range("F16").Value = 659
